$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update sheet name
$ws.Name = "Through 2022-03-20"

# Update the "March (through 03-19)" label to "March (through 03-20)"
$ws.Range("A4").Value = "March (through 03-20)"

# Update March row (row 4) values for columns C..I (2016..2022)
$ws.Range("C4").Value = 29
$ws.Range("D4").Value = 35
$ws.Range("E4").Value = 40
$ws.Range("F4").Value = 21
$ws.Range("G4").Value = 42
$ws.Range("H4").Value = 54
$ws.Range("I4").Value = 88

# Update Total row (row 5) values for columns C..I (2016..2022)
$ws.Range("C5").Value = 116
$ws.Range("D5").Value = 166
$ws.Range("E5").Value = 177
$ws.Range("F5").Value = 100
$ws.Range("G5").Value = 183
$ws.Range("H5").Value = 396
$ws.Range("I5").Value = 388
